# Auto-generated edit script: updates LeveProfit calculation columns (H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4529.7334
$ws.Range("I64").Value = 3722.25
$ws.Range("J64").Value = 4653.9614
$ws.Range("K64").Value = 3722.25
$ws.Range("L64").Value = 4653.9614
$ws.Range("M64").Value = -3474.25
$ws.Range("N64").Value = -5149.9614
$ws.Range("H67").Value = 4529.7334
$ws.Range("I67").Value = 3722.25
$ws.Range("J67").Value = 4653.9614
$ws.Range("K67").Value = 3722.25
$ws.Range("L67").Value = 4653.9614
$ws.Range("M67").Value = -2864.25
$ws.Range("N67").Value = -6369.9614
$ws.Range("H76").Value = 4400
$ws.Range("I76").Value = 4800
$ws.Range("J76").Value = 4342.857
$ws.Range("K76").Value = 4800
$ws.Range("L76").Value = 4342.857
$ws.Range("M76").Value = -4485
$ws.Range("N76").Value = -4972.857
$ws.Range("H79").Value = 4400
$ws.Range("I79").Value = 4800
$ws.Range("J79").Value = 4342.857
$ws.Range("K79").Value = 4800
$ws.Range("L79").Value = 4342.857
$ws.Range("M79").Value = -3708
$ws.Range("N79").Value = -6526.857
$ws.Range("H112").Value = 6693
$ws.Range("J112").Value = 7452.9565
$ws.Range("L112").Value = 22358.8695
$ws.Range("N112").Value = -24574.8695
$ws.Range("H137").Value = 2067.2334
$ws.Range("I137").Value = 1839.6522
$ws.Range("J137").Value = 2815
$ws.Range("K137").Value = 5518.9566
$ws.Range("L137").Value = 8445
$ws.Range("M137").Value = -2968.9566
$ws.Range("N137").Value = -13545
$ws.Range("H138").Value = 2227623.8
$ws.Range("I138").Value = 5409231
$ws.Range("J138").Value = 6501.717
$ws.Range("K138").Value = 16227693
$ws.Range("L138").Value = 19505.151
$ws.Range("M138").Value = -16222553
$ws.Range("N138").Value = -29785.151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3001.1667
$ws.Range("J88").Value = 3001.1667
$ws.Range("L88").Value = 3001.1667
$ws.Range("N88").Value = -3813.1667
$ws.Range("H91").Value = 3001.1667
$ws.Range("J91").Value = 3001.1667
$ws.Range("L91").Value = 3001.1667
$ws.Range("N91").Value = -5809.1667
$ws.Range("H132").Value = 2017.1154
$ws.Range("I132").Value = 1712.9025
$ws.Range("J132").Value = 3151
$ws.Range("K132").Value = 5138.7075
$ws.Range("L132").Value = 9453
$ws.Range("M132").Value = -2608.7075
$ws.Range("N132").Value = -14513

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2388.2307
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 3144.717
$ws.Range("I134").Value = 1888.5405
$ws.Range("J134").Value = 6049.625
$ws.Range("K134").Value = 5665.6215
$ws.Range("L134").Value = 18148.875
$ws.Range("M134").Value = -3130.6215
$ws.Range("N134").Value = -23218.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3553.8206
$ws.Range("I31").Value = 2038.7916
$ws.Range("J31").Value = 5977.8667
$ws.Range("K31").Value = 2038.7916
$ws.Range("L31").Value = 5977.8667
$ws.Range("M31").Value = -1743.7916
$ws.Range("N31").Value = -6567.8667
$ws.Range("H34").Value = 3553.8206
$ws.Range("I34").Value = 2038.7916
$ws.Range("J34").Value = 5977.8667
$ws.Range("K34").Value = 2038.7916
$ws.Range("L34").Value = 5977.8667
$ws.Range("M34").Value = -1836.7916
$ws.Range("N34").Value = -6381.8667
$ws.Range("H62").Value = 261002.5
$ws.Range("I62").Value = 500005
$ws.Range("J62").Value = 22000
$ws.Range("K62").Value = 500005
$ws.Range("L62").Value = 22000
$ws.Range("M62").Value = -499381
$ws.Range("N62").Value = -23248
$ws.Range("H65").Value = 261002.5
$ws.Range("I65").Value = 500005
$ws.Range("J65").Value = 22000
$ws.Range("K65").Value = 2500025
$ws.Range("L65").Value = 110000
$ws.Range("M65").Value = -2496905
$ws.Range("N65").Value = -116240
$ws.Range("H134").Value = 1633.5862
$ws.Range("I134").Value = 1442.2916
$ws.Range("K134").Value = 4326.8748
$ws.Range("M134").Value = -1791.8748

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4000
$ws.Range("I56").Value = 4000
$ws.Range("K56").Value = 4000
$ws.Range("M56").Value = -3470
$ws.Range("H98").Value = 2353.9
$ws.Range("I98").Value = 5045
$ws.Range("J98").Value = 1681.125
$ws.Range("K98").Value = 15135
$ws.Range("L98").Value = 5043.375
$ws.Range("M98").Value = -13637
$ws.Range("N98").Value = -8039.375
$ws.Range("H113").Value = 3258.4255
$ws.Range("I113").Value = 884.8333
$ws.Range("J113").Value = 3605.7805
$ws.Range("K113").Value = 2654.4999
$ws.Range("L113").Value = 10817.3415
$ws.Range("M113").Value = -484.4998999999998
$ws.Range("N113").Value = -15157.3415

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6019.0356
$ws.Range("I70").Value = 4624.923
$ws.Range("K70").Value = 4624.923
$ws.Range("M70").Value = -4354.923
$ws.Range("H73").Value = 6019.0356
$ws.Range("I73").Value = 4624.923
$ws.Range("K73").Value = 4624.923
$ws.Range("M73").Value = -3688.923
$ws.Range("H80").Value = 3070.7222
$ws.Range("I80").Value = 2957.3333
$ws.Range("J80").Value = 3297.5
$ws.Range("K80").Value = 2957.3333
$ws.Range("L80").Value = 3297.5
$ws.Range("M80").Value = -1959.3333
$ws.Range("N80").Value = -5293.5
$ws.Range("H83").Value = 3070.7222
$ws.Range("I83").Value = 2957.3333
$ws.Range("J83").Value = 3297.5
$ws.Range("K83").Value = 14786.6665
$ws.Range("L83").Value = 16487.5
$ws.Range("M83").Value = -9794.666499999999
$ws.Range("N83").Value = -26471.5
$ws.Range("H102").Value = 1965.5714
$ws.Range("I102").Value = 1913.8334
$ws.Range("K102").Value = 1913.8334
$ws.Range("M102").Value = -291.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 20000
$ws.Range("J3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("N3").Value = -20224
$ws.Range("H4").Value = 50000
$ws.Range("I4").Value = 50000
$ws.Range("K4").Value = 50000
$ws.Range("M4").Value = -49887
$ws.Range("H10").Value = 21126.334
$ws.Range("I10").Value = 380
$ws.Range("J10").Value = 31499.5
$ws.Range("K10").Value = 380
$ws.Range("L10").Value = 31499.5
$ws.Range("N10").Value = -31779.5
$ws.Range("M10").Value = -240
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20340
$ws.Range("H20").Value = 30025000
$ws.Range("I20").Value = 30025000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 30025000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -30024774
$ws.Range("N20").ClearContents()
$ws.Range("H21").Value = 19000
$ws.Range("I21").Value = 19000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -18826
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 255
$ws.Range("I22").Value = 255
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 255
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 40
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 255
$ws.Range("I27").Value = 255
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 255
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -148
$ws.Range("N27").ClearContents()
$ws.Range("H28").Value = 50000
$ws.Range("I28").Value = 50000
$ws.Range("K28").Value = 50000
$ws.Range("M28").Value = -49768
$ws.Range("H37").Value = 50000
$ws.Range("I37").Value = 50000
$ws.Range("K37").Value = 50000
$ws.Range("M37").Value = -49893
$ws.Range("H132").Value = 7860.5264
$ws.Range("I132").Value = 8587.5
$ws.Range("J132").Value = 3983.3333
$ws.Range("K132").Value = 25762.5
$ws.Range("L132").Value = 11949.9999
$ws.Range("M132").Value = -23232.5
$ws.Range("N132").Value = -17009.9999

Write-Host "Applied profit refresh updates"